# Update the "Joystick Mapping" sheet for v3.5 (Mar 13, 2019):
#  - bump version label
#  - mark Flip Manipulator with an asterisk
#  - relabel the "Manipulator Vertical" starting-state cell as "Starting Configuration"
#  - shift the Ball Intake legend entries up one row, replacing "High/Low Gear"
#  - fill the now-empty legend cells with "-"
#  - clear the stray "Milky Manipulator" label
#  - split "Drive to Vision Target" into Bumper Cam / Manipulator Cam rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H1").Value = "v3.5 (Mar 13, 2019)"

$ws.Range("B3").Value = "Flip Manipulator*"
$ws.Range("B5").Value = "Starting Configuration"

$ws.Range("E16").Value = "Ball Intake Unfold/Fold"
$ws.Range("E17").Value = "Ball Intake Wheels"
$ws.Range("E18").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "-"

$ws.Range("B25").Value = ""

$ws.Range("H29").Value = "Drive to Vision Target - Bumper Cam"
$ws.Range("H28").Value = "Drive to Vision Target - Manipulator Cam"
